{"js": "// The document's placeholder tokens (e.g. \"${numero_informe}\") were split\n// across more runs by Word's editor (adding spell-check proofErr markers\n// around words like \"N\u00ba\" and the snake_case placeholder names), but the\n// only actual *content* change in this revision is in the big notification\n// paragraph: the literal \"${numero_resolucion}\" placeholder was replaced\n// with the real resolution number \"0222-2022-CU-UNH\", and the trailing\n// space before the final period was dropped.\n//\n//   \"...aprobado con resoluci\u00f3n N\u00b0 ${numero_resolucion}. \"\n//   -> \"...aprobado con resoluci\u00f3n N\u00b0 0222-2022-CU-UNH.\"\n\nconst body = context.document.body;\n\nconst results = body.search(\"${numero_resolucion}. \", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  results.items[0].insertText(\"0222-2022-CU-UNH.\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# The document's placeholder tokens (e.g. \"${numero_informe}\") were split\n# across more runs by Word's editor (adding spell-check proofErr markers\n# around words like \"N\u00ba\" and the snake_case placeholder names), but the\n# only actual *content* change in this revision is in the big notification\n# paragraph: the literal \"${numero_resolucion}\" placeholder was replaced\n# with the real resolution number \"0222-2022-CU-UNH\", and the trailing\n# space before the final period was dropped.\n#\n#   \"...aprobado con resoluci\u00f3n N\u00b0 ${numero_resolucion}. \"\n#   -> \"...aprobado con resoluci\u00f3n N\u00b0 0222-2022-CU-UNH.\"\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.Text = \"`$`{numero_resolucion`}. \"\n$find.Replacement.Text = \"0222-2022-CU-UNH.\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n"}
